{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that holds the (bookmark-split) \"Baz changes\" text.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Baz changes\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Could not find the 'Baz changes' paragraph\");\n}\n\n// Insert the \"Version Management\" heading paragraph right before it.\ntarget.insertParagraph(\"Version Management \", \"Before\");\n\n// Insert the version-management body paragraph right before it too.\nconst bodyText =\n  \"Version management is managing the various components such as code files, \" +\n  \"configuration files, data file, documentation, media items or any digital item \" +\n  \"that is required to build an entire system. \";\nconst bodyPara = target.insertParagraph(bodyText, \"Before\");\nawait context.sync();\n\n// Bold the leading \"Version management\" phrase of the new paragraph.\nconst matches = bodyPara.search(\"Version management\", { matchCase: true });\nmatches.load(\"items\");\nawait context.sync();\nif (matches.items.length > 0) {\n  matches.items[0].font.bold = true;\n}\nawait context.sync();\n\n// Collapse the \"Baz chan\" / bookmark / \"ges\" runs into a single run,\n// re-homing the _GoBack bookmark at the start of the paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\ntarget.insertText(\"Baz changes\", \"Replace\");\nawait context.sync();\n\ntarget.getRange(\"Start\").insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Get-BazChangesRange($doc) {\n    $found = $null\n    foreach ($p in $doc.Paragraphs) {\n        if ($p.Range.Text -eq \"Baz changes`r\") {\n            $found = $p.Range\n        }\n    }\n    return $found\n}\n\nfunction Get-ParagraphIndexAt($doc, $startPos) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        if ($doc.Paragraphs($i).Range.Start -eq $startPos) {\n            return $i\n        }\n    }\n    return -1\n}\n\n# Locate the paragraph that holds the (bookmark-split) \"Baz changes\" text.\n$target = Get-BazChangesRange $d\nif ($target -eq $null) {\n    throw \"Could not find the 'Baz changes' paragraph\"\n}\n\n# --- Insert the \"Version Management\" heading paragraph right before it ---\n$target.InsertParagraphBefore()\n$target = Get-BazChangesRange $d\n$idx = Get-ParagraphIndexAt $d $target.Start\n$heading = $d.Paragraphs($idx - 1).Range\n$heading.Text = \"Version Management \"\n\n# --- Insert the version-management body paragraph right before it ---\n$target = Get-BazChangesRange $d\n$target.InsertParagraphBefore()\n$target = Get-BazChangesRange $d\n$idx = Get-ParagraphIndexAt $d $target.Start\n$bodyPara = $d.Paragraphs($idx - 1).Range\n$bodyText = \"Version management is managing the various components such as code files, configuration files, data file, documentation, media items or any digital item that is required to build an entire system. \"\n$bodyPara.Text = $bodyText\n\n# Bold just the leading \"Version management\" phrase of the new paragraph.\n$boldLen = (\"Version management\").Length\n$boldRange = $d.Range($bodyPara.Start, $bodyPara.Start + $boldLen)\n$boldRange.Font.Bold = 1\n\n# --- Collapse the \"Baz chan\" / bookmark / \"ges\" runs into one run ---\n$target = Get-BazChangesRange $d\n$startPos = $target.Start\n$endPos = $target.End\n\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$editRange = $d.Range($startPos, $endPos)\n$editRange.Text = \"Baz changes\"\n\n$bmRange = $d.Range($startPos, $startPos)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n"}
